$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.707523345947266
$ws.Range("B1").Value = 2.832193613052368
$ws.Range("C1").Value = 3.114429950714111
$ws.Range("D1").Value = 3.511049270629883
$ws.Range("E1").Value = 1.604800820350647
